# Applies the changes described by the commit:
#  - DIMS sheet: TYPE column values "STR32_ID" -> "ID"; narrow the TYPE column.
#  - DVDND_TYP_ENUM sheet: insert a new "Reset"/"-" /"Reset value" row right
#    after the header, widen the NAME column, and keep the autofilter /
#    filter-database range in sync with the now-larger table.

$wb = $excel.ActiveWorkbook

# --- DIMS sheet -----------------------------------------------------------
$dims = $wb.Worksheets.Item("DIMS")

foreach ($r in 2,4,6,8,10,12,14) {
    $dims.Cells.Item($r, 3).Value = "ID"
}

# Narrow column C (TYPE) from 13 to ~7.8 characters.
$dims.Columns.Item(3).ColumnWidth = 7.0

# --- DVDND_TYP_ENUM sheet ---------------------------------------------------
$dvdnd = $wb.Worksheets.Item("DVDND_TYP_ENUM")

# Insert a new row right below the header and copy the formatting that the
# following data row already has, so the new row matches the table style.
$dvdnd.Rows.Item(2).Insert()
$dvdnd.Range("A3:C3").Copy()
$dvdnd.Range("A2:C2").PasteSpecial(-4122)
$dvdnd.Rows.Item(2).RowHeight = 36

$dvdnd.Range("A2").Value = "-"
$dvdnd.Range("B2").Value = "Reset"
$dvdnd.Range("C2").Value = "Reset value"

# Widen column B (NAME) from ~7.8 to ~9.1 characters.
$dvdnd.Columns.Item(2).ColumnWidth = 8.333333333333332

# Re-apply the autofilter over the now 16-row table (A1:C16).
$dvdnd.AutoFilterMode = $false
$dvdnd.Range("A1:C16").AutoFilter()

# Keep the workbook-level hidden _FilterDatabase name for this sheet in sync
# with the new autofilter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "DVDND_TYP_ENUM!_FilterDatabase") {
        $n.RefersTo = "='DVDND_TYP_ENUM'!`$A`$1:`$C`$16"
    }
}

Write-Host "Edit complete"
